# Weekly refresh: prepend a new observation row at the top of the data
# table (row 2), pushing every existing data row down by one. The table
# is sorted newest-first, so this mirrors a new week's record being added
# to the consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stash the current row 2 (it will become row 3 once we make room).
$oldRow2 = @()
for ($c = 1; $c -le 18; $c++) {
    $oldRow2 += ,($ws.Cells.Item(2, $c).Value())
}

# Insert a blank row above the current row 3. Excel copies the new row's
# formatting from the row directly above it (row 2, an ordinary data row
# carrying the column-D date style) instead of the bold header row, so no
# stray styles get introduced.
$ws.Range("A3:R3").Insert()

# Move the old row 2 data down into the freshly made row 3.
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(3, $c).Value = $oldRow2[$c - 1]
}

# Overwrite row 2 with this week's new record.
$ws.Range("A2").Value = 8
$ws.Range("B2").Value = "Terminal La Palmera de La Serena"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44922
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 100112030
$ws.Range("G2").Value = "Poroto granado"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 520
$ws.Range("K2").Value = 29000
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = 29500
$ws.Range("N2").Value = "$/malla 25 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 1180
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
